$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the CLAVE column so rows 3-7 share the same key as row 2 (EMP27707)
$ws.Range("A3").Value = "EMP27707"
$ws.Range("A4").Value = "EMP27707"
$ws.Range("A5").Value = "EMP27707"
$ws.Range("A6").Value = "EMP27707"
$ws.Range("A7").Value = "EMP27707"

# Update the selection to match the edited range
$ws.Range("A3:A7").Select()
